$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for the Price/Volume columns so values stay as text
# (matching the source data which stores figures as plain strings), rather than
# being auto-converted to numbers/percentages by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# Updated crypto price / 1h volume figures (GitHub Actions data refresh)
$ws.Range("D2").Value = "308.66"
$ws.Range("E2").Value = "-2.56%"
$ws.Range("D3").Value = "37.81"
$ws.Range("E3").Value = "-4.15%"
$ws.Range("D4").Value = "5.093"
$ws.Range("E4").Value = "-0.36%"
$ws.Range("D5").Value = "0.07898"
$ws.Range("E5").Value = "-3.64%"
$ws.Range("D6").Value = "1.971"
$ws.Range("E6").Value = "-3.09%"
$ws.Range("D7").Value = "4.349"
$ws.Range("E7").Value = "1.63%"
$ws.Range("D8").Value = "8.249"
$ws.Range("E8").Value = "-0.23%"
$ws.Range("D9").Value = "3.198"
$ws.Range("E9").Value = "-2.95%"
$ws.Range("D10").Value = "0.9331"
$ws.Range("E10").Value = "0.00%"
$ws.Range("D11").Value = "0.1302"
$ws.Range("D12").Value = "0.1949"
$ws.Range("E12").Value = "-2.08%"
$ws.Range("D13").Value = "0.08831"
$ws.Range("E13").Value = "-3.16%"
$ws.Range("D14").Value = "0.03423"
$ws.Range("E14").Value = "-3.77%"
$ws.Range("D15").Value = "0.09729"
$ws.Range("E15").Value = "-0.84%"
$ws.Range("D16").Value = "0.001396"
$ws.Range("E16").Value = "-0.26%"
$ws.Range("D17").Value = "0.005950"
$ws.Range("E17").Value = "-5.62%"
$ws.Range("E18").Value = "1,777.56%"
$ws.Range("D19").Value = "3.598"
$ws.Range("E19").Value = "-1.67%"
$ws.Range("D20").Value = "0.3438"
$ws.Range("E20").Value = "-0.63%"
$ws.Range("D21").Value = "0.1295"
$ws.Range("E21").Value = "-0.62%"
$ws.Range("D22").Value = "4.999"
$ws.Range("E22").Value = "1.92%"
$ws.Range("D23").Value = "0.2486"
$ws.Range("E23").Value = "1.47%"
$ws.Range("D24").Value = "0.04304"
$ws.Range("E24").Value = "-0.59%"
$ws.Range("D25").Value = "0.001217"
$ws.Range("E25").Value = "-0.62%"
$ws.Range("D26").Value = "0.004614"
$ws.Range("E26").Value = "-3.59%"
$ws.Range("E27").Value = "176.20%"
$ws.Range("E39").Value = "3.43%"
$ws.Range("D40").Value = "0.05059"
$ws.Range("E40").Value = "-3.55%"
$ws.Range("D41").Value = "0.007511"
$ws.Range("E41").Value = "-0.11%"
$ws.Range("D42").Value = "0.009896"
$ws.Range("E42").Value = "1.20%"
$ws.Range("D43").Value = "0.1358"
$ws.Range("E43").Value = "-1.45%"
$ws.Range("E44").Value = "-2.93%"
$ws.Range("D45").Value = "0.007991"
$ws.Range("E45").Value = "-15.90%"
$ws.Range("D46").Value = "0.00006559"
$ws.Range("E46").Value = "1.76%"
$ws.Range("E47").Value = "-0.04%"
$ws.Range("E48").Value = "8.23%"
$ws.Range("E50").Value = "-0.04%"
$ws.Range("E51").Value = "-0.04%"
